$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = $ws.Range("A2").Value2
